$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C23 keeps its text "slope estimate" (shared string index changes internally but text is same) - no-op, but set explicitly for safety
$ws.Range("C23").Value = "slope estimate"

# Add new text to C24: "which plot to assess? & which seems reasonable? : linearity "
$ws.Range("C24").Value = "which plot to assess? & which seems reasonable? : linearity "

# Update C21: "interpret slope, interpret intercept" -> "interpret slope, interpret intercept, assumptions of lm check"
$ws.Range("C21").Value = "interpret slope, interpret intercept, assumptions of lm check"

# Update the active selection from C17 to B16
$ws.Range("B16").Select()
